# Generate Report for Handoff
# - Status text moves from "Handed back: in sync with en-US" to "Ready for handoff"
# - Timestamp cells bump forward to the new handoff-generation time
# - The wide "status/date" columns shrink to the narrower handoff-report width

$wb = $excel.ActiveWorkbook

# Target stored column width is 17.2159881591797 characters. Excel's
# ColumnWidth setter quantizes through a pixel grid, so the character
# value that lands closest to that target (17.16666...) is fed in here.
$newWidth = 16.333333333333336

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-18 02:52:27"
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-18 02:52:22"
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-18 02:52:27"
$dede.Columns.Item(3).ColumnWidth = $newWidth
